# Regenerate s_vals data to filter save games.
# Updates columns B (TB), C (d2S), D (K), E (IP) and the recalculated
# G (sum) column for data rows 2-11. Column A (date) and F (Win) are
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ B = 0.3048080303191223;  C = 0.3127903958511391;  D = 0.1575252929769615; E = 0.496779210170732;  G = 1.271902929317955 }
    3  = @{ B = 3.230985683306322;   C = 1.667794583268128;   D = 0.8054896365839992; E = 0.496779210170732;  G = 6.201049113329182 }
    4  = @{ B = 0.6753301551942219;  C = 1.667794583268128;   D = 3.900430680208489;  E = 8.660232485948974;  G = 14.90378790461981 }
    5  = @{ B = 1.459612070389937;   C = 1.667794583268128;   D = 0.1575252929769615; E = 0.496779210170732;  G = 3.781711156805759 }
    6  = @{ B = 1.459612070389937;   C = 1.667794583268128;   D = 0.1575252929769615; E = 0.496779210170732;  G = 3.781711156805759 }
    7  = @{ B = 1.459612070389937;   C = 1.667794583268128;   D = 0.1575252929769615; E = 0.496779210170732;  G = 3.781711156805759 }
    8  = @{ B = 1.459612070389937;   C = 0.3127903958511391;  D = 0.8054896365839992; E = 0.496779210170732;  G = 3.074671312995807 }
    9  = @{ B = 0.6753301551942219;  C = 0.04240448674262143; D = 0.8054896365839992; E = 0.496779210170732;  G = 2.020003488691574 }
    10 = @{ B = 1.459612070389937;   C = 0.3127903958511391;  D = 0.1575252929769615; E = 0.496779210170732;  G = 2.42670696938877 }
    11 = @{ B = 1.459612070389937;   C = 0.3127903958511391;  D = 3.900430680208489;  E = 0.496779210170732;  G = 6.169612356620297 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}
